$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renewal-case data fix: BI_SYMBOL/PD_SYMBOL/UM_SYMBOL/MP_SYMBOL on row 5
# move from code "N" to code "K".
$ws.Range("AE5:AH5").Value = "K"

# Sheet view update: scroll the window so column Z is the left-most visible
# column, and move the active selection to the single cell AI2.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 26
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AI2").Select()
